$d = $word.ActiveDocument

# Locate the "Git basic terms" paragraph robustly via Find.
$findRng = $d.Content
$found = $findRng.Find.Execute("Git basic terms", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Git basic terms' paragraph"
}

# The paragraph holding the found text is where the new "Repository" /
# "commit" paragraphs get appended. We also fold in the trailing (blank)
# paragraph of the document body, since it collapses away once the new
# text paragraphs take its place as the final paragraph.
$paraIndex = $findRng.Paragraphs.Item(1).Index
$startPara = $d.Paragraphs.Item($paraIndex)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$rng = $d.Range($startPara.Range.Start, $lastPara.Range.End)

# Rewrite that span as three paragraphs:
#   1) "Git basic terms" (unchanged text, same paraId, bookmark removed)
#   2) "Repository" (new)
#   3) "commit" (new) carrying the _GoBack bookmark that used to sit on
#      paragraph 1, now moved to the end of the new final paragraph.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="123E4F71"><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>Git basic terms</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>Repository</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>commit</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$rng.InsertXML($xml)
